$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1340
$ws.Range("E2").Value = 88
$ws.Range("F2").Value = 88
$ws.Range("G2").Value = 55
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 1794
$ws.Range("L2").Value = 624
$ws.Range("M2").Value = 1171
$ws.Range("N2").Value = 1137
$ws.Range("O2").Value = 34
$ws.Range("P2").Value = 172
$ws.Range("Q2").Value = 147
$ws.Range("R2").Value = 83
$ws.Range("S2").Value = -71
$ws.Range("T2").Value = 47
$ws.Range("U2").Value = 100
$ws.Range("V2").Value = 354
$ws.Range("W2").Value = 6.57
$ws.Range("X2").Value = 2.42
$ws.Range("Y2").Value = 2.68
$ws.Range("Z2").Value = 1.78
$ws.Range("AA2").Value = 53.26
$ws.Range("AB2").Value = 595.31
$ws.Range("AC2").Value = 95
$ws.Range("AD2").Value = 72.33
$ws.Range("AE2").Value = 3964
$ws.Range("AF2").Value = 1.74
$ws.Range("AG2").Value = 75
$ws.Range("AH2").Value = 1.09
$ws.Range("AI2").Value = 71.56
$ws.Range("AJ2").Value = 31496785

# Row 3
$ws.Range("D3").Value = 1055
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = 14
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1717
$ws.Range("L3").Value = 562
$ws.Range("M3").Value = 1155
$ws.Range("N3").Value = 1124
$ws.Range("O3").Value = 31
$ws.Range("P3").Value = 172
$ws.Range("Q3").Value = -32
$ws.Range("R3").Value = -211
$ws.Range("S3").Value = -13
$ws.Range("T3").Value = 146
$ws.Range("U3").Value = -178
$ws.Range("V3").Value = 369
$ws.Range("W3").Value = 1.47
$ws.Range("X3").Value = 1.35
$ws.Range("Y3").Value = 1.25
$ws.Range("Z3").Value = 0.81
$ws.Range("AA3").Value = 48.64
$ws.Range("AB3").Value = 588.23
$ws.Range("AC3").Value = 45
$ws.Range("AD3").Value = 142.5
$ws.Range("AE3").Value = 3919
$ws.Range("AF3").Value = 1.63
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 0.78
$ws.Range("AI3").Value = 101.35
$ws.Range("AJ3").Value = 31496785

# Row 4
$ws.Range("D4").Value = 1097
$ws.Range("E4").Value = 55
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 17
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 1796
$ws.Range("L4").Value = 644
$ws.Range("M4").Value = 1152
$ws.Range("N4").Value = 1119
$ws.Range("O4").Value = 33
$ws.Range("P4").Value = 172
$ws.Range("Q4").Value = 115
$ws.Range("R4").Value = -98
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 131
$ws.Range("U4").Value = -16
$ws.Range("V4").Value = 395
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 1.55
$ws.Range("Y4").Value = 0.9
$ws.Range("Z4").Value = 0.97
$ws.Range("AA4").Value = 55.92
$ws.Range("AB4").Value = 585.93
$ws.Range("AC4").Value = 32
$ws.Range("AD4").Value = 147.68
$ws.Range("AE4").Value = 3902
$ws.Range("AF4").Value = 1.21
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 1.06
$ws.Range("AI4").Value = 142.86
$ws.Range("AJ4").Value = 31496785

# Row 5
$ws.Range("D5").Value = 1522
$ws.Range("E5").Value = 115
$ws.Range("F5").Value = 115
$ws.Range("G5").Value = 270
$ws.Range("H5").Value = 194
$ws.Range("I5").Value = 183
$ws.Range("J5").Value = 11
$ws.Range("K5").Value = 2016
$ws.Range("L5").Value = 680
$ws.Range("M5").Value = 1335
$ws.Range("N5").Value = 1293
$ws.Range("O5").Value = 43
$ws.Range("P5").Value = 172
$ws.Range("Q5").Value = 62
$ws.Range("R5").Value = 278
$ws.Range("S5").Value = -125
$ws.Range("T5").Value = 64
$ws.Range("U5").Value = -2
$ws.Range("V5").Value = 275
$ws.Range("W5").Value = 7.57
$ws.Range("X5").Value = 12.78
$ws.Range("Y5").Value = 15.19
$ws.Range("Z5").Value = 10.2
$ws.Range("AA5").Value = 50.96
$ws.Range("AB5").Value = 687.12
$ws.Range("AC5").Value = 582
$ws.Range("AD5").Value = 8.82
$ws.Range("AE5").Value = 4508
$ws.Range("AF5").Value = 1.14
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.95
$ws.Range("AI5").Value = 15.65
$ws.Range("AJ5").Value = 31496785

# Row 6
$ws.Range("D6").Value = 1960
$ws.Range("E6").Value = 215
$ws.Range("F6").Value = 215
$ws.Range("G6").Value = 206
$ws.Range("H6").Value = 147
$ws.Range("I6").Value = 104
$ws.Range("K6").Value = 2105
$ws.Range("L6").Value = 664
$ws.Range("M6").Value = 1441
$ws.Range("N6").Value = 1355
$ws.Range("P6").Value = 172
$ws.Range("Q6").Value = 92
$ws.Range("R6").Value = -231
$ws.Range("S6").Value = 13
$ws.Range("T6").Value = 143
$ws.Range("U6").Value = -51
$ws.Range("V6").Value = 317
$ws.Range("W6").Value = 10.98
$ws.Range("X6").Value = 7.49
$ws.Range("Y6").Value = 7.86
$ws.Range("Z6").Value = 7.13
$ws.Range("AA6").Value = 46.09
$ws.Range("AB6").Value = 725.79
$ws.Range("AC6").Value = 330
$ws.Range("AD6").Value = 10.58
$ws.Range("AE6").Value = 4727
$ws.Range("AF6").Value = 0.74
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 1.43
$ws.Range("AI6").Value = 13.78
$ws.Range("AJ6").Value = 31496785

# Row 7: clear all data columns D:AJ
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns D:AJ
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns D:AJ
$ws.Range("D9:AJ9").ClearContents()
